# Corrected some selection scopes
# The underlying per-vintage error vectors (columns B:K on rows 3-22) were
# shifted down by one row (a new, freshly computed vintage is inserted at
# row 3, and the former row-22 tail value becomes a brand new row 23),
# while the row labels in column A stay attached to their original row
# number. A new date label is also introduced for the new trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift the numeric data (columns B:K) down by one row, working from the
#    bottom up so we never clobber a source row before it has been read.
for ($r = 22; $r -ge 3; $r--) {
    $srcRow = $ws.Range("B" + $r + ":K" + $r)
    $vals = $srcRow.Value2

    $dstRow = $ws.Range("B" + ($r + 1) + ":K" + ($r + 1))
    $dstRow.ClearContents()

    for ($c = 2; $c -le 11; $c++) {
        $v = $vals.GetValue(1, $c - 1)
        if ($v -ne $null) {
            $ws.Cells.Item($r + 1, $c).Value2 = $v
        }
    }
}

# 2. Populate the newly freed row 3 with the freshly computed values.
$ws.Range("B3:K3").ClearContents()
$ws.Range("B3").Value2 = -3.02587959343237
$ws.Range("C3").Value2 = -1.776533676533679
$ws.Range("D3").Value2 = 2.832862923122467
$ws.Range("E3").Value2 = 2.815775584009379
$ws.Range("F3").Value2 = -1.250511408339877
$ws.Range("G3").Value2 = 0.4135197463245859
$ws.Range("H3").Value2 = -0.1376344086021527

# 3. Label the new trailing row (row 23) the same way the other vintage rows
#    are labelled (bold, centered, bordered), with its own new date string.
#    Copy the format from the row above first, then overwrite the value so
#    we end up re-using the existing shared cell style instead of minting a
#    new, unused one.
$ws.Range("A22").Copy($ws.Range("A23"))
$ws.Range("A23").Value = "2020-05-15 00:00:00_diff"
